$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-28 down to 16-29.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(15, 3).Value = "La Araucanía"
$ws.Cells.Item(15, 4).Value = 44484
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = 300000000
$ws.Cells.Item(15, 7).Value = "Espárragos"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 1200
$ws.Cells.Item(15, 12).Value = 1200
$ws.Cells.Item(15, 13).Value = 1200
$ws.Cells.Item(15, 14).Value = "$/kilo"
$ws.Cells.Item(15, 15).Value = "Región del Maule"
$ws.Cells.Item(15, 16).Value = 1200
$ws.Cells.Item(15, 17).Value = 1
$ws.Cells.Item(15, 18).Value = "Hortaliza"
